$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheets
$ws1.Name = "Tier system_v01-00"
$ws2.Name = "Tier system_v01-00_SUMMARY"

# Update cell text values
$ws1.Range("F3").Value = "FISH, Immuno Fluorescence, Expression of FP"
$ws1.Range("C7").Value = "Pioneer"
$ws1.Range("E7").Value = "Development of novel unproved technology or of new gold-standard; full reproducibility ofmicroscopy set up and image acquisition settings"

# Remove comment on C7
$ws1.Range("C7").Comment.Delete()

# Update selections on sheets
$ws1.Range("H5").Select()
$ws2.Range("I5").Select()

# Make sheet1 active (so tabSelected moves there, activeTab removed/0)
$ws1.Activate()

Write-Host "done"
